# Updated symbol list with GitHub Actions crypto scrape refresh.
# D column holds prices as plain text (e.g. "244.08"); a leading
# apostrophe forces Excel to keep the numeric-looking text as a
# string instead of silently re-typing it as a Number, then the
# style is reset to Normal so no stray "quote prefix" formatting
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.408"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05972"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.466"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.513"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8137"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9149"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.07430"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03223"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "'0.09366"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.850"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001569"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04678"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.006133"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.004999"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0009819"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.00007802"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.615"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.1302"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002394"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.03929"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003001"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007191"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "'0.00005237"
$ws.Range("D45").Style = "Normal"
$ws.Range("D49").Value = "'0.002289"
$ws.Range("D49").Style = "Normal"
